$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This sheet is a chronological bank-statement ledger: row 1 always holds the
# newest movement, and every time a new movement is recorded, the existing
# rows are pushed down by one. Reproduce that: shift rows 1..45 down to
# 2..46 (columns A:G - the raw data columns; column H holds a derived
# formula that only ever lives on row 1), then populate the new row 1 with
# the new movement, and fix up the trailing blank row / dimension / H column.
# ---------------------------------------------------------------------------

# 1) Make sure a (still blank) row 46 exists with the same formatting as the
#    other trailing blank rows (41-45), by copying row 45's formatting down.
$ws.Range("A45").Copy() | Out-Null
$ws.Range("A46").PasteSpecial(-4122) | Out-Null

# 2) Shift every data row down by one: work bottom-up so we never overwrite
#    a source row before it has been copied.
for ($r = 45; $r -ge 1; $r--) {
    $src = $ws.Range("A" + $r + ":G" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":G" + ($r + 1))
    $dst.Value2 = $src.Value2
}

# 3) Column H only ever carries the CONCATENATE helper formula, and only on
#    row 1. After the shift, drop the stray shared-formula remnant that
#    landed on H2 (old H1/H2 formulas do not move down with the data).
$ws.Range("H2").ClearContents() | Out-Null

# 4) Write the new top movement into row 1.
$nbsp = [char]0x00A0
$ws.Range("A1").Value2 = 41712
$ws.Range("B1").Value2 = "2200555126/0995935959"
$ws.Range("C1").Value2 = "D"
$ws.Range("D1").Value2 = "0007745172"
$ws.Range("E1").Value2 = "AG. NORTE"
$ws.Range("F1").Value2 = "6.00$nbsp$nbsp"
$ws.Range("G1").Value2 = "1995.56"

# 5) Re-assert the helper formula on H1 (values shifted, so it recomputes
#    against the new row 1 inputs).
$ws.Range("H1").Formula = "=CONCATENATE(""array('mo_fecha' => new \DateTime('"",TEXT(A1,""yyyy-mm-dd""),""'), 'mo_concepto' => '"",B1,""', 'mo_tipo' => '"",C1,""', 'mo_documento' => '"",D1,""', 'mo_oficina' => '"",E1,""', 'mo_monto' => "",TRIM(F1),"", 'mo_saldo' => "",G1,"", 'mo_fecha_crea' => new \DateTime('"",TEXT(NOW(),""yyyy-mm-dd H:m:s""),""'), 'mo_quien_crea' => 1, 'mo_fecha_modifica' => NULL, 'mo_quien_modifica' => NULL, 'mo_borrado_logico' => false),"")"

# 6) Selection now targets just H1 (the shared range used to span H1:H2).
$ws.Range("H1").Select() | Out-Null
